# The scraper re-ran and appended two more match rows (4 and 5) to the
# "Josh Philippe †" sheet. The new rows duplicate the data already present
# in rows 2 and 3 respectively (same venue/date/result/opponent/stats), so
# the cleanest way to reproduce them with matching cell types/formatting
# (several columns hold numeric-looking values like "0"/"3"/"0.00" that
# must remain text, not become numbers) is to copy the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:K2").Copy($ws.Range("A4:K4"))
$ws.Range("A3:K3").Copy($ws.Range("A5:K5"))
